# Add the new "dvd" row to the active sheet ("pasty rezystywne"),
# extending the used range from A1:C2 to A1:C3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "dvd"
$ws.Range("B3").Value = 32
$ws.Range("C3").Value = 12
